$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 65.016
$ws.Range("D2").Value = 65.016
$ws.Range("E2").Value = 2.65640079
$ws.Range("F2").Value = 0.0006583399999999999
$ws.Range("G2").Value = 0.04264630000000001
$ws.Range("H2").Value = 2.8224482
$ws.Range("I2").Value = 9.558719219284198
$ws.Range("J2").Value = 9.558719219284198
$ws.Range("K2").Value = 0.395088426762209
$ws.Range("L2").Value = 0.00007814186170389435
$ws.Range("M2").Value = 0.007085029834641331
$ws.Range("N2").Value = 0.8247088289406709

$ws.Range("C3").Value = 111.411
$ws.Range("D3").Value = 111.411
$ws.Range("E3").Value = 1.56493938
$ws.Range("F3").Value = 0.0004032900000000001
$ws.Range("G3").Value = 0.04468881
$ws.Range("H3").Value = 5.112078170000001
$ws.Range("I3").Value = 20.24381117504535
$ws.Range("J3").Value = 20.24381117504535
$ws.Range("K3").Value = 0.2746311276630839
$ws.Range("L3").Value = 0.00004751982779886764
$ws.Range("M3").Value = 0.008309330734852306
$ws.Range("N3").Value = 1.851909330208276

$ws.Range("C4").Value = 42.202
$ws.Range("D4").Value = 84.374
$ws.Range("E4").Value = 2.11203137
$ws.Range("F4").Value = 0.00134037
$ws.Range("G4").Value = 0.02806220000000001
$ws.Range("H4").Value = 1.23351562
$ws.Range("I4").Value = 9.355999718401595
$ws.Range("J4").Value = 18.69913853047198
$ws.Range("K4").Value = 0.5171385935415701
$ws.Range("L4").Value = 0.0001359707731505236
$ws.Range("M4").Value = 0.005912061944570307
$ws.Range("N4").Value = 0.5107474579666893

$ws.Range("C5").Value = 64.185
$ws.Range("D5").Value = 125.068
$ws.Range("E5").Value = 1.39943506
$ws.Range("F5").Value = 0.00084256
$ws.Range("G5").Value = 0.02686643
$ws.Range("H5").Value = 1.78094903
$ws.Range("I5").Value = 12.74563557560886
$ws.Range("J5").Value = 23.19653973021527
$ws.Range("K5").Value = 0.2700096485585489
$ws.Range("L5").Value = 0.0001091410012697252
$ws.Range("M5").Value = 0.005545948403907018
$ws.Range("N5").Value = 0.6976701909553732

$ws.Range("C6").Value = 25.863
$ws.Range("D6").Value = 103.36
$ws.Range("E6").Value = 1.76239613
$ws.Range("F6").Value = 0.0027586
$ws.Range("G6").Value = 0.0177257
$ws.Range("H6").Value = 0.48852123
$ws.Range("I6").Value = 6.918243497672882
$ws.Range("J6").Value = 27.65862251449626
$ws.Range("K6").Value = 0.5116546991245751
$ws.Range("L6").Value = 0.0002923369646502147
$ws.Range("M6").Value = 0.004716807987144228
$ws.Range("N6").Value = 0.2575916055772217

$ws.Range("C7").Value = 34.198
$ws.Range("D7").Value = 125.683
$ws.Range("E7").Value = 1.4016614
$ws.Range("F7").Value = 0.00177757
$ws.Range("G7").Value = 0.01510953
$ws.Range("H7").Value = 0.53868211
$ws.Range("I7").Value = 7.516468139578142
$ws.Range("J7").Value = 25.61934546123031
$ws.Range("K7").Value = 0.2915120878691205
$ws.Range("L7").Value = 0.0002290896706034062
$ws.Range("M7").Value = 0.003527867561095318
$ws.Range("N7").Value = 0.2393767752610564

$ws.Range("C8").Value = 17.736
$ws.Range("D8").Value = 106.237
$ws.Range("E8").Value = 1.76099322
$ws.Range("F8").Value = 0.00433683
$ws.Range("G8").Value = 0.0127027
$ws.Range("H8").Value = 0.24554356
$ws.Range("I8").Value = 5.563924987026808
$ws.Range("J8").Value = 33.34160304324754
$ws.Range("K8").Value = 0.5950308720627353
$ws.Range("L8").Value = 0.000511768656971349
$ws.Range("M8").Value = 0.003893556175232315
$ws.Range("N8").Value = 0.1530169136507125

$ws.Range("C9").Value = 22.335
$ws.Range("D9").Value = 114.255
$ws.Range("E9").Value = 1.54139866
$ws.Range("F9").Value = 0.00273074
$ws.Range("G9").Value = 0.01005951
$ws.Range("H9").Value = 0.23449934
$ws.Range("I9").Value = 5.201906846035868
$ws.Range("J9").Value = 23.48747829454317
$ws.Range("K9").Value = 0.3161308627061897
$ws.Range("L9").Value = 0.0003877987109832113
$ws.Range("M9").Value = 0.002342197461540946
$ws.Range("N9").Value = 0.1120015352352791

$ws.Range("C10").Value = 12.944
$ws.Range("D10").Value = 103.422
$ws.Range("E10").Value = 1.84404614
$ws.Range("F10").Value = 0.005592629999999999
$ws.Range("G10").Value = 0.008965979999999998
$ws.Range("H10").Value = 0.12878856
$ws.Range("I10").Value = 4.454065861962832
$ws.Range("J10").Value = 35.63913669077948
$ws.Range("K10").Value = 0.672208549587134
$ws.Range("L10").Value = 0.0005716640472190664
$ws.Range("M10").Value = 0.002968582308193681
$ws.Range("N10").Value = 0.08826241086951397

$ws.Range("C11").Value = 16.31
$ws.Range("D11").Value = 102.1
$ws.Range("E11").Value = 1.72702627
$ws.Range("F11").Value = 0.0034409
$ws.Range("G11").Value = 0.006950700000000001
$ws.Range("H11").Value = 0.11836115
$ws.Range("I11").Value = 3.784735687498431
$ws.Range("J11").Value = 21.30647708951355
$ws.Range("K11").Value = 0.3563532894379963
$ws.Range("L11").Value = 0.0004999090908346227
$ws.Range("M11").Value = 0.001647232978812192
$ws.Range("N11").Value = 0.05358529320835826

$ws.Range("C12").Value = 10.098
$ws.Range("D12").Value = 100.786
$ws.Range("E12").Value = 1.91720592
$ws.Range("F12").Value = 0.00719583
$ws.Range("G12").Value = 0.007192019999999999
$ws.Range("H12").Value = 0.08168604
$ws.Range("I12").Value = 3.729920683916524
$ws.Range("J12").Value = 37.29752417130806
$ws.Range("K12").Value = 0.7325737882595844
$ws.Range("L12").Value = 0.0007832487406243879
$ws.Range("M12").Value = 0.002532487356693378
$ws.Range("N12").Value = 0.06092564508109

$ws.Range("C13").Value = 13.192
$ws.Range("D13").Value = 91.884
$ws.Range("E13").Value = 1.92569394
$ws.Range("F13").Value = 0.00409328
$ws.Range("G13").Value = 0.005331300000000001
$ws.Range("H13").Value = 0.07385673
$ws.Range("I13").Value = 3.369048305751416
$ws.Range("J13").Value = 20.13559003322192
$ws.Range("K13").Value = 0.4130035035071562
$ws.Range("L13").Value = 0.0007117063773548524
$ws.Range("M13").Value = 0.001385306964317073
$ws.Range("N13").Value = 0.03724574927217798
